# 2.0 Sep update 3
# Three forecast-model rows (the old row 4 "Holt Winter"/2, old row 7
# "ARIMA"/5, and old row 9 "Prophet"/7 entries) were dropped from the
# comparison table. Deleting the corresponding sheet rows shifts the
# remaining data up so the sheet goes from 9 rows to 6 rows, matching
# the refreshed model-comparison export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so row indices of not-yet-deleted rows stay valid.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(4).Delete()

# Leave the same selection state captured in the saved workbook (entire row 7).
$ws.Rows.Item(7).Select() | Out-Null
